$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 81; this shifts all existing rows 81..108 down to 82..109
$ws.Rows(81).Insert()

# Populate the newly inserted row 81 with the new record
$ws.Cells.Item(81, 1).Value = 5
$ws.Cells.Item(81, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(81, 3).Value = "Maule"
$ws.Cells.Item(81, 4).Value = 45146
$ws.Cells.Item(81, 5).Value = 7
$ws.Cells.Item(81, 6).Value = 100112040
$ws.Cells.Item(81, 7).Value = "Cilantro"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 200
$ws.Cells.Item(81, 11).Value = 8000
$ws.Cells.Item(81, 12).Value = 8000
$ws.Cells.Item(81, 13).Value = 8000
$ws.Cells.Item(81, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(81, 15).Value = "Región Metropolitana"
$ws.Cells.Item(81, 16).Value = 222
$ws.Cells.Item(81, 17).Value = 36
$ws.Cells.Item(81, 18).Value = "Hortaliza"
